$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 1.960204196686623
$ws.Range("B4").Value = 1.633503497238852
$ws.Range("B5").Value = 0.3267006994477705
$ws.Range("B9").Value = 0.2347046349630462
$ws.Range("B10").Value = 0.001116790388228642
$ws.Range("B11").Value = 0.004472629178150887
$ws.Range("B13").Value = 0.157375484050397
$ws.Range("B14").Value = 0.07045138194024789
$ws.Range("B15").Value = 0.07732858432060842
$ws.Range("B16").Value = 0.07546342384598012
$ws.Range("B17").Value = 0.03770864205288811
$ws.Range("B18").Value = 0.003350371164685926
